# Updated cryptos list on Fri Jun 28 23:12:59 UTC 2024 with GitHub Actions
# Refreshes price/volume figures in columns D and E, and fixes the
# ranking order of a couple of coin rows (Aptos/Monero, Stacks/OKB swap).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that can look numeric (e.g. "41.28"); force
# text formatting first so Excel doesn't silently coerce it to a number
# and drop things like trailing zeros.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.379.35'
$ws.Range("E2").Value = '  -1.99%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.372.08'
$ws.Range("E3").Value = '  -2.23%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '567.09'
$ws.Range("E5").Value = '  -2.22%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.04'
$ws.Range("E6").Value = '  -6.60%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.372.96'
$ws.Range("E8").Value = '  -2.24%  '
$ws.Range("E9").Value = '  -0.78%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.50'
$ws.Range("E10").Value = '  -3.66%  '
$ws.Range("E11").Value = '  -3.69%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.386'
$ws.Range("E12").Value = '  -1.43%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.946.98'
$ws.Range("E13").Value = '  -2.31%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.04'
$ws.Range("E14").Value = '  +0.31%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.123'
$ws.Range("E15").Value = '  +0.91%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.375.06'
$ws.Range("E16").Value = '  -2.49%  '
$ws.Range("E17").Value = '  -3.77%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '60.498.07'
$ws.Range("E18").Value = '  -2.04%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.17'
$ws.Range("E19").Value = '  -1.93%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.77'
$ws.Range("E20").Value = '  -4.39%  '
$ws.Range("E21").Value = '  -5.54%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '385.99'
$ws.Range("E22").Value = '  -1.04%  '
$ws.Range("E23").Value = '  -2.26%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '73.04'
$ws.Range("E24").Value = '  -0.15%  '
$ws.Range("E25").Value = '  -0.14%  '
$ws.Range("E26").Value = '  -8.26%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.518.14'
$ws.Range("E27").Value = '  -1.97%  '
$ws.Range("E28").Value = '  -1.34%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.09%  '
$ws.Range("E30").Value = '  -5.18%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.90'
$ws.Range("E31").Value = '  -4.24%  '
$ws.Range("E32").Value = '  -2.38%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.40'
$ws.Range("E33").Value = '  -9.39%  '
$ws.Range("E34").Value = '  -0.01%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.44'
$ws.Range("E35").Value = '  -2.72%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.402.19'
$ws.Range("E36").Value = '  -2.03%  '
$ws.Range("B37").Value = 'Monero'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '168.10'
$ws.Range("E37").Value = '  +0.82%  '
$ws.Range("B38").Value = 'Aptos'
$ws.Range("C38").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.86'
$ws.Range("E38").Value = '  -2.40%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.92'
$ws.Range("E39").Value = '  -5.68%  '
$ws.Range("E40").Value = '  -5.02%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0767'
$ws.Range("E41").Value = '  -2.66%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '26.99'
$ws.Range("E42").Value = '  -0.15%  '
$ws.Range("E43").Value = '  -0.05%  '
$ws.Range("E44").Value = '  -2.08%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.41'
$ws.Range("E45").Value = '  -2.30%  '
$ws.Range("B46").Value = 'OKB'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '41.28'
$ws.Range("E46").Value = '  -2.44%  '
$ws.Range("B47").Value = 'Stacks'
$ws.Range("C47").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.69'
$ws.Range("E47").Value = '  -1.57%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.511.85'
$ws.Range("E48").Value = '  -3.59%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.12'
$ws.Range("E49").Value = '  -4.51%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '23.10'
$ws.Range("E50").Value = '  -0.72%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.73'
$ws.Range("E51").Value = '  -3.21%  '
